$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.517.89"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").Value = "3.898.48"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'602.60"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").Value = "'168.78"
$ws.Range("E6").Value = "  +1.65%  "
$ws.Range("D7").Value = "3.899.60"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.530"
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("E10").Value = "  -0.20%  "
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("D12").Value = "'0.459"
$ws.Range("E12").Value = "  -0.26%  "
$ws.Range("E13").Value = "  -1.31%  "
$ws.Range("D14").Value = "'37.11"
$ws.Range("E14").Value = "  -0.45%  "
$ws.Range("D15").Value = "4.552.43"
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("D16").Value = "3.918.07"
$ws.Range("E16").Value = "  +0.16%  "
$ws.Range("D17").Value = "68.465.60"
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D18").Value = "'18.15"
$ws.Range("E18").Value = "  +5.65%  "
$ws.Range("D19").Value = "'7.43"
$ws.Range("E19").Value = "  -0.62%  "
$ws.Range("E20").Value = "  +0.26%  "
$ws.Range("D21").Value = "'10.89"
$ws.Range("E21").Value = "  -1.20%  "
$ws.Range("D22").Value = "'472.91"
$ws.Range("E22").Value = "  -2.88%  "
$ws.Range("D23").Value = "'0.740"
$ws.Range("E23").Value = "  +2.34%  "
$ws.Range("D24").Value = "'0.0000167"
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").Value = "'83.77"
$ws.Range("E25").Value = "  -0.90%  "
$ws.Range("E26").Value = "  +1.07%  "
$ws.Range("D27").Value = "'12.28"
$ws.Range("E27").Value = "  +1.78%  "
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "'10.00"
$ws.Range("E29").Value = "  -1.07%  "
$ws.Range("E30").Value = "  +1.15%  "
$ws.Range("D31").Value = "4.049.67"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").Value = "'7.87"
$ws.Range("E32").Value = "  +1.87%  "
$ws.Range("D33").Value = "'31.53"
$ws.Range("E33").Value = "  -0.92%  "
$ws.Range("E34").Value = "  -2.51%  "
$ws.Range("D35").Value = "'9.44"
$ws.Range("E35").Value = "  +1.89%  "
$ws.Range("D36").Value = "3.873.31"
$ws.Range("E36").Value = "  +0.51%  "
$ws.Range("E37").Value = "  -1.65%  "
$ws.Range("D38").Value = "'3.66"
$ws.Range("E38").Value = "  +14.36%  "
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("E40").Value = "  +2.50%  "
$ws.Range("D41").Value = "'5.93"
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("E42").Value = "  +0.16%  "
$ws.Range("D43").Value = "'0.314"
$ws.Range("E43").Value = "  -0.37%  "
$ws.Range("D44").Value = "'430.07"
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("D45").Value = "'2.00"
$ws.Range("E45").Value = "  +0.97%  "
$ws.Range("D46").Value = "'0.000298"
$ws.Range("E46").Value = "  +12.90%  "
$ws.Range("D48").Value = "'8.62"
$ws.Range("E48").Value = "  +1.27%  "
$ws.Range("D49").Value = "'47.26"
$ws.Range("E49").Value = "  -2.14%  "
$ws.Range("D50").Value = "'27.05"
$ws.Range("E50").Value = "  +5.89%  "
$ws.Range("D51").Value = "'143.65"
$ws.Range("E51").Value = "  +0.96%  "
